# "Generate Report for Handback"
#
# The localization-status report previously showed every file as still
# "Ready for handoff" with no handback info recorded. This run populates
# the per-language sheets with the handback results: status flips to
# "Handed back: in sync with en-US", the (previously empty) "Latest Target
# File" / "Latest Handback File" columns are filled in (mirroring the
# source .md / handoff .xlf names, each as a live hyperlink like the other
# filename columns), and "Latest Handback DateTime" gets a real timestamp
# instead of the zero-date placeholder.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Per-sheet "Latest Handback DateTime" stamp.
$handbackTimes = @{
    "zh-cn" = "2016-03-19 06:36:14"
    "de-de" = "2016-03-19 06:36:19"
}

foreach ($sheetName in $handbackTimes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $handbackTime = $handbackTimes[$sheetName]

    for ($row = 2; $row -le 3; $row++) {
        $aCell = $ws.Range("A$row")
        $dCell = $ws.Range("D$row")
        $fCell = $ws.Range("F$row")
        $gCell = $ws.Range("G$row")
        $cCell = $ws.Range("C$row")
        $hCell = $ws.Range("H$row")

        # Find the existing hyperlinks on the source (A) and handoff-target
        # (D) cells so the new Target/Handback columns can reuse the same
        # address + display text.
        $aLink = $null
        $dLink = $null
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range.Address() -eq $aCell.Address()) { $aLink = $h }
            if ($h.Range.Address() -eq $dCell.Address()) { $dLink = $h }
        }

        # Status -> handed back.
        $cCell.Value = $statusText

        # Latest Target File (F) mirrors the source file name/link (A).
        $ws.Hyperlinks.Add($fCell, $aLink.Address(), "", "", $aLink.TextToDisplay())

        # Latest Handback File (G) mirrors the handoff target file name/link (D).
        $ws.Hyperlinks.Add($gCell, $dLink.Address(), "", "", $dLink.TextToDisplay())

        # Latest Handback DateTime (H) - plain text timestamp, same shape
        # as the existing "Latest Handoff Datetime" column.
        $hCell.Value = $handbackTime
    }
}
